$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the data (columns A, B, E, F, G, H, Q, R) among rows 2, 3 and 4:
#   new Row2 = old Row4
#   new Row3 = old Row2
#   new Row4 = old Row3
# All other columns are identical across these three rows, so only these
# columns need to be updated.

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

# Capture the original values before overwriting anything.
$orig2 = @{}
$orig3 = @{}
$orig4 = @{}
foreach ($col in $cols) {
    $orig2[$col] = $ws.Range("$col`2").Value2
    $orig3[$col] = $ws.Range("$col`3").Value2
    $orig4[$col] = $ws.Range("$col`4").Value2
}

foreach ($col in $cols) {
    $ws.Range("$col`2").Value2 = $orig4[$col]
    $ws.Range("$col`3").Value2 = $orig2[$col]
    $ws.Range("$col`4").Value2 = $orig3[$col]
}
